$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "40.789.28"
$ws.Range("E2").Value = "  +3.63%  "

Set-TextValue $ws.Range("D3") "2.214.97"
$ws.Range("E3").Value = "  +2.63%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue $ws.Range("D5") "229.74"
$ws.Range("E5").Value = "  +0.50%  "

Set-TextValue $ws.Range("D6") "0.633"
$ws.Range("E6").Value = "  +1.91%  "

Set-TextValue $ws.Range("D7") "64.54"
$ws.Range("E7").Value = "  +0.55%  "

$ws.Range("E8").Value = "  +0.11%  "

Set-TextValue $ws.Range("D9") "0.405"
$ws.Range("E9").Value = "  +1.93%  "

Set-TextValue $ws.Range("D10") "0.0868"
$ws.Range("E10").Value = "  +1.21%  "

$ws.Range("E11").Value = "  +0.15%  "

Set-TextValue $ws.Range("D12") "2.544.01"
$ws.Range("E12").Value = "  +2.54%  "

Set-TextValue $ws.Range("D13") "15.92"
$ws.Range("E13").Value = "  +0.07%  "

Set-TextValue $ws.Range("D14") "22.25"
$ws.Range("E14").Value = "  +0.03%  "

Set-TextValue $ws.Range("D15") "0.822"
$ws.Range("E15").Value = "  +0.86%  "

Set-TextValue $ws.Range("D16") "5.62"

Set-TextValue $ws.Range("D17") "2.221.87"
$ws.Range("E17").Value = "  +3.33%  "

Set-TextValue $ws.Range("D18") "40.671.84"
$ws.Range("E18").Value = "  +3.49%  "

Set-TextValue $ws.Range("D19") "74.15"

Set-TextValue $ws.Range("D20") "0.0₃0903"
$ws.Range("E20").Value = "  +6.11%  "

Set-TextValue $ws.Range("D21") "6.16"
$ws.Range("E21").Value = "  +0.41%  "

Set-TextValue $ws.Range("D22") "251.44"
$ws.Range("E22").Value = "  +8.78%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D24") "2.37"
$ws.Range("E24").Value = "  -5.56%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D25") "2.38"
$ws.Range("E25").Value = "  +1.00%  "

Set-TextValue $ws.Range("D26") "9.72"
$ws.Range("E26").Value = "  +0.92%  "

Set-TextValue $ws.Range("D27") "173.12"

$ws.Range("E28").Value = "  +3.20%  "

Set-TextValue $ws.Range("D29") "20.37"
$ws.Range("E29").Value = "  +2.30%  "

$ws.Range("E30").Value = "  +2.40%  "

Set-TextValue $ws.Range("D31") "2.82"
$ws.Range("E31").Value = "  +4.38%  "

$ws.Range("E32").Value = "  +1.42%  "

Set-TextValue $ws.Range("D33") "4.67"
$ws.Range("E33").Value = "  +1.12%  "

Set-TextValue $ws.Range("D34") "7.16"
$ws.Range("E34").Value = "  +1.14%  "

Set-TextValue $ws.Range("D35") "4.78"
$ws.Range("E35").Value = "  -0.09%  "

Set-TextValue $ws.Range("D36") "0.0632"
$ws.Range("E36").Value = "  +2.29%  "

$ws.Range("E37").Value = "  +6.45%  "

Set-TextValue $ws.Range("D38") "2.46"
$ws.Range("E38").Value = "  +1.54%  "

Set-TextValue $ws.Range("D40") "4.92"
$ws.Range("E40").Value = "  +14.30%  "

$ws.Range("E41").Value = "  +1.60%  "

Set-TextValue $ws.Range("D42") "8.62"
$ws.Range("E42").Value = "  +10.87%  "

Set-TextValue $ws.Range("D43") "101.43"
$ws.Range("E43").Value = "  -2.39%  "

$ws.Range("E44").Value = "  +4.25%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D45") "17.34"
$ws.Range("E45").Value = "  -2.65%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D46") "1.511.28"
$ws.Range("E46").Value = "  -1.80%  "

$ws.Range("E47").Value = "  +1.54%  "

Set-TextValue $ws.Range("D48") "1.11"
$ws.Range("E48").Value = "  +1.86%  "

$ws.Range("E49").Value = "  +0.16%  "

Set-TextValue $ws.Range("D50") "0.000205"
$ws.Range("E50").Value = "  +38.30%  "

Set-TextValue $ws.Range("D51") "9.64"
$ws.Range("E51").Value = "  +12.25%  "
